$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a "Category" header in A1, matching the style of the other header cells.
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# The category labels in A2:A46 should not carry the bold/bordered header
# style anymore - reset them to the default (unstyled) look.
$ws.Range("A2:A46").Style = "Normal"
